# Updates the workbook to:
#   1. Set the "Förändrad" (column C) date value to 45186 for every data row.
#   2. Add the article id (column A) as the friendly-name second argument of
#      every HYPERLINK() formula found in columns S-Y.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (data starts on row 2; row 1
# holds the headers, row 0 is an empty spacer row present in this file).
$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {

    # --- 1) bump the "Förändrad" date value in column C -------------------
    $ws.Cells.Item($r, 3).Value = 45186

    # --- 2) patch every HYPERLINK formula in columns S (19) .. Y (25) -----
    $beteckning = $ws.Cells.Item($r, 1).Value()

    for ($c = 19; $c -le 25; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula) {
            $formula = $cell.Formula
            if ($formula -like '*HYPERLINK(*' -and $formula -notlike '*,*') {
                $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $beteckning + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
